# Insert a new data row at row 293 (pushing the existing rows 293-352 down
# to 294-353) and populate it with the new observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(293).Insert()

$ws.Cells.Item(293, 1).Value = 4
$ws.Cells.Item(293, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(293, 3).Value = "Los Lagos"
$ws.Cells.Item(293, 4).Value = 44637
$ws.Cells.Item(293, 5).Value = 10
$ws.Cells.Item(293, 6).Value = 100114001
$ws.Cells.Item(293, 7).Value = "Papa"
$ws.Cells.Item(293, 8).Value = "Patagonia"
$ws.Cells.Item(293, 9).Value = "1a (cosecha)"
$ws.Cells.Item(293, 10).Value = 300
$ws.Cells.Item(293, 11).Value = 7000
$ws.Cells.Item(293, 12).Value = 8000
$ws.Cells.Item(293, 13).Value = 7500
$ws.Cells.Item(293, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(293, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(293, 16).Value = 300
$ws.Cells.Item(293, 17).Value = 25
$ws.Cells.Item(293, 18).Value = "Hortaliza"
